$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.610.80"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "2.603.20"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.67"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.36"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  -5.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.77"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.11"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "3.072.10"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  -7.70%  "
$ws.Range("D16").Value = "63.441.55"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "2.617.21"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -4.13%  "
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.46"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.22"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "575.68"
$ws.Range("E27").Value = "  +9.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.81"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.86"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.30"
$ws.Range("E42").Value = "  -2.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "155.73"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.79"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.626"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.89"
$ws.Range("E51").Value = "  -3.98%  "
